$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update version/date/publisher metadata, replace the
# duplicated "Contact" row with a "Jurisdiction" row, and drop the now
# redundant second "Contact" row entirely (21 rows -> 20 rows). ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Cells.Item(3, 2).Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now "Alvearie Team"
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes
# "Jurisdiction" / "United States of America"
$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# Row 11 was the duplicate "Contact" / "No display for ContactDetail" row -
# remove it entirely, shifting all following rows up by one.
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements": the top-level Extension row's Short/Definition
# columns (K2/L2) are updated. ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Cells.Item(2, 11).Value = "Longterm Care Duration"
$elements.Cells.Item(2, 12).Value = "Number of weeks for which the employee is eligible for long-term disability (LTD) benefits"
